$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "loc"
$ws.Range("C2").Value = "loc.png"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0

$ws.Range("I2").Select()
